$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.037.57'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '1.555.13'
$ws.Range('E3').Value = '  +1.23%  '
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $origStyle
$ws.Range('E4').Value = '  +0.21%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.04'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('E7').Value = '  +0.42%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.84'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('E10').Value = '  +0.98%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0857'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '1.774.75'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '1.549.99'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('E15').Value = '  +2.27%  '
$ws.Range('D16').Value = '26.945.28'
$ws.Range('E16').Value = '  +0.68%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.71'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +1.29%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '217.43'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('D19').Value = '0.0₃0689'
$ws.Range('E19').Value = '  +1.20%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.24'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  +0.54%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.05'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  +1.14%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.11'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +1.90%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.59'
$ws.Range('D26').Style = $origStyle
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.88'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('E30').Value = '  +3.13%  '
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '1.440.67'
$ws.Range('E33').Value = '  +5.77%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.04'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('E35').Value = '  +3.98%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.965'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  +1.17%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('E42').Value = '  -0.42%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.986'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('E44').Value = '  +3.40%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.91'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +1.91%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.75'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').Value = '1.689.11'
$ws.Range('E47').Value = '  +1.15%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.22'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('E49').Value = '  +3.84%  '
$ws.Range('D50').Value = '0.0₆0100'
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('E51').Value = '  +1.66%  '
